$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reservation dates in column D (rows 3-8)
$ws.Range("D4").Value = "08/01/2020"
$ws.Range("D3").Value = "08/02/2020"
$ws.Range("D5").Value = "07/31/2020"
$ws.Range("D6").Value = "07/31/2020"
$ws.Range("D7").Value = "07/31/2020"
$ws.Range("D8").Value = "07/31/2020"

# Update the active selection to D3
$ws.Range("D3").Select()
